{"js": "// Insert the new \"Further, we developed...\" sentences in the cover-letter\n// paragraph, right before \"We hope that this tool allows specialists and \".\n// The inserted text inherits the formatting of the surrounding run\n// (Times New Roman, en-US) because Word.InsertLocation.before splits off a\n// new run that copies the formatting of the range it is anchored to.\n\nconst anchorText = \"We hope that this tool allows specialists and \";\nconst insertedText =\n  \"Further, we developed a new metric to quantify Darwinian shortfalls. \" +\n  \"We illustrate all the functionalities \" +\n  \"of our package by constructing phylogenies for the four most speciose freshwater ecoregions of the world, besides to map the Darwinian shortfalls for all the assemblages in freshwater basins in the world. \";\n\nconst body = context.document.body;\nconst results = body.search(anchorText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Anchor text \"' + anchorText + '\" not found in document.');\n}\n\n// The anchor text is unique in this letter (appears once), so just use the\n// first (only) hit.\nconst target = results.items[0];\ntarget.insertText(insertedText, Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "# Insert the new \"Further, we developed...\" sentences in the cover-letter\n# paragraph, right before \"We hope that this tool allows specialists and \".\n\n$d = $word.ActiveDocument\n\n$insertedText = \"Further, we developed a new metric to quantify Darwinian shortfalls. \" `\n  + \"We illustrate all the functionalities \" `\n  + \"of our package by constructing phylogenies for the four most speciose freshwater ecoregions of the world, besides to map the Darwinian shortfalls for all the assemblages in freshwater basins in the world. \"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWildcards = $false\n$found = $rng.Find.Execute(\"We hope that this tool allows specialists and \")\n\nif (-not $found) {\n  throw \"Anchor text not found in document.\"\n}\n\n# $rng now collapses to exactly the matched text (\"We hope ... and \");\n# inserting before it places the new sentences right ahead of it, inheriting\n# the surrounding run's formatting (Times New Roman, en-US).\n$rng.InsertBefore($insertedText)\n"}
